# Update plots for each sample
# ------------------------------------------------------------------
# This script reproduces the data-level edits captured in the commit:
#   - peak_table:   the "m_height" measurement for CYP2D6_49 (row 4) drops
#                    to 800 (re-measured peak height)
#   - allele_table: the previously-undetected CYP2D6_49 / Reverse / T
#                    (wildtype) peak (row 6) is now detected -> fills in
#                    peak/size/height/status and clears the failure message
#   - marker_table: the genotype call for CYP2D6_49 (row 4) is now resolved
#                    to "TT" / "wildtype"
#   - genotype_result: the overall sample genotype is now "*1/*1"
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$peak = $wb.Sheets("peak_table")
$allele = $wb.Sheets("allele_table")
$marker = $wb.Sheets("marker_table")
$genotype = $wb.Sheets("genotype_result")

# ---- peak_table: CYP2D6_49 (row 4) wildtype peak height measurement ----
$peak.Range("N4").Value = 800

# ---- allele_table: row 6 (CYP2D6_003 / CYP2D6_49 / Reverse / T / wildtype) ----
# Peak height re-measured lower (K = min_height column used for detection)
$allele.Range("K6").Value = 800
# Peak is now flagged as detected
$allele.Range("M6").Value = $true
# Detected peak/size/height values
$allele.Range("N6").Value = 17
$allele.Range("O6").Value = 38.87
$allele.Range("P6").Value = 872
# Status now "ok"
$allele.Range("Q6").Value = "ok"
# No failure message anymore
$allele.Range("R6").Value = ""

# ---- marker_table: row 4 (CYP2D6_003 / CYP2D6_49) genotype call resolved ----
$marker.Range("G4").Value = "TT"
$marker.Range("H4").Value = "wildtype"

# ---- genotype_result: overall sample genotype ----
$genotype.Range("B2").Value = "*1/*1"

# Leave the selection on the cell that was actually re-measured, mirroring
# where the author was working when the peak plot got updated.
$peak.Activate() | Out-Null
$peak.Range("N4").Select() | Out-Null
